# Apply updated crypto price/volume figures to columns D (Price) and E (Volume(1h)).
# Values that look like plain decimal numbers (single '.') are prefixed with a literal
# apostrophe so Excel stores them as text, matching the workbook's original inline-string
# cells instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.685.17'
$ws.Range("E2").Value = '  -1.51%  '
$ws.Range("D3").Value = '2.906.32'
$ws.Range("E3").Value = '  -2.19%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''528.73'
$ws.Range("E5").Value = '  -2.84%  '
$ws.Range("D6").Value = '''143.42'
$ws.Range("E6").Value = '  -5.68%  '
$ws.Range("D8").Value = '''0.554'
$ws.Range("E8").Value = '  -3.42%  '
$ws.Range("D9").Value = '2.913.48'
$ws.Range("E9").Value = '  -2.30%  '
$ws.Range("E10").Value = '  -5.06%  '
$ws.Range("D11").Value = '''6.03'
$ws.Range("E11").Value = '  -1.92%  '
$ws.Range("D12").Value = '''0.360'
$ws.Range("E12").Value = '  -2.91%  '
$ws.Range("D13").Value = '3.416.49'
$ws.Range("E13").Value = '  -2.17%  '
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("D15").Value = '60.645.91'
$ws.Range("E15").Value = '  -1.71%  '
$ws.Range("D16").Value = '''22.78'
$ws.Range("E16").Value = '  -4.20%  '
$ws.Range("D17").Value = '2.916.92'
$ws.Range("E17").Value = '  -2.07%  '
$ws.Range("E18").Value = '  -4.04%  '
$ws.Range("D19").Value = '''5.03'
$ws.Range("E19").Value = '  -2.76%  '
$ws.Range("D20").Value = '''11.71'
$ws.Range("E20").Value = '  -2.74%  '
$ws.Range("D21").Value = '''361.16'
$ws.Range("E21").Value = '  -5.56%  '
$ws.Range("D22").Value = '''6.63'
$ws.Range("E22").Value = '  -1.25%  '
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").Value = '''5.68'
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").Value = '''64.89'
$ws.Range("E25").Value = '  -1.52%  '
$ws.Range("E26").Value = '  -3.32%  '
$ws.Range("E27").Value = '  -4.03%  '
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("D29").Value = '''7.88'
$ws.Range("E29").Value = '  -5.44%  '
$ws.Range("D30").Value = '0.0₃0848'
$ws.Range("E30").Value = '  -10.10%  '
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").Value = '''1.70'
$ws.Range("E32").Value = '  -1.69%  '
$ws.Range("D33").Value = '''19.82'
$ws.Range("E33").Value = '  -3.30%  '
$ws.Range("D34").Value = '''148.57'
$ws.Range("E34").Value = '  -7.66%  '
$ws.Range("D35").Value = '''4.35'
$ws.Range("E35").Value = '  -6.69%  '
$ws.Range("E36").Value = '  -6.33%  '
$ws.Range("E37").Value = '  -6.70%  '
$ws.Range("D38").Value = '''1.20'
$ws.Range("E38").Value = '  -5.37%  '
$ws.Range("D39").Value = '''37.91'
$ws.Range("E39").Value = '  +1.57%  '
$ws.Range("E40").Value = '  -4.63%  '
$ws.Range("D41").Value = '''3.72'
$ws.Range("E41").Value = '  -5.28%  '
$ws.Range("D42").Value = '2.297.37'
$ws.Range("E42").Value = '  -4.82%  '
$ws.Range("E43").Value = '  -2.48%  '
$ws.Range("D44").Value = '''0.0587'
$ws.Range("E44").Value = '  -1.41%  '
$ws.Range("D45").Value = '''20.51'
$ws.Range("E45").Value = '  -7.96%  '
$ws.Range("D46").Value = '''0.997'
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '''4.99'
$ws.Range("E47").Value = '  -1.81%  '
$ws.Range("E48").Value = '  -3.87%  '
$ws.Range("D49").Value = '''10.33'
$ws.Range("E49").Value = '  -1.21%  '
$ws.Range("D50").Value = '''0.0922'
$ws.Range("E50").Value = '  -3.20%  '
$ws.Range("D51").Value = '''249.75'
$ws.Range("E51").Value = '  -7.64%  '
